{"js": "// Fix the Section 508 sufficiency check: the template text was missing the\n// `.section508Sufficient` property access, which caused the generated\n// document to always evaluate the ternary as truthy (\"Yes\").\n//\n// Before: {sensitiveInformation.section508 ? `Yes` : `No, complete ...`}\n// After:  {sensitiveInformation.section508.section508Sufficient ? `Yes` : `No, complete ...`}\n\nconst body = context.document.body;\n\n// Find the run containing the old property-access expression. Search on the\n// unique substring right before the insertion point so we don't depend on\n// matching the (very long) rest of the sentence.\nconst searchResults = body.search(\"sensitiveInformation.section508\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find 'sensitiveInformation.section508' text to update.\");\n}\n\n// Insert the missing property access immediately after the matched text,\n// right before the \" ? `Yes` : ...\" ternary.\nconst target = searchResults.items[0];\ntarget.insertText(\".section508Sufficient\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Fix the Section 508 sufficiency check: the template text was missing the\n# `.section508Sufficient` property access, which caused the generated\n# document to always evaluate the ternary as truthy (\"Yes\").\n#\n# Before: {sensitiveInformation.section508 ? `Yes` : `No, complete ...`}\n# After:  {sensitiveInformation.section508.section508Sufficient ? `Yes` : `No, complete ...`}\n\n$d = $word.ActiveDocument\n\n# Locate the expression text (search on a short, unique substring rather than\n# the whole long sentence) and position the range right after it.\n$range = $d.Content\n$found = $range.Find.Execute(\"sensitiveInformation.section508\")\n\nif ($found) {\n    # Collapse the found range to its end point, then insert the missing\n    # property access immediately before the \" ? `Yes` : ...\" ternary.\n    $range.Collapse(0)  # wdCollapseEnd\n    $range.InsertAfter(\".section508Sufficient\")\n}\n"}
